$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# D-column values are written as text (using a Text number format that is reset
# back to the default "Normal" style afterwards) so that values which look like
# plain numbers (e.g. "9.30", "0.0621") are not silently coerced into doubles and
# keep their original formatting (trailing zeros, etc.), matching values like
# "36.918.65" that already behave as text because Excel cannot parse them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.918.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.981.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

$ws.Range("E6").Value = "  +1.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.93"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.63%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("E11").Value = "  +0.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.90%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.268.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("E16").Value = "  +3.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.981.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.799.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.46%  "

$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  +4.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("E30").Value = "  +17.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.80%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("E39").Value = "  -7.57%  "

$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.372.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "

$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.21%  "
